$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.5569716666666666
$ws.Range("H2").Value = 1.670915
$ws.Range("I2").Value = 0.3425729542218473
$ws.Range("J2").Value = 0.3425729542218473
$ws.Range("M2").Value = 3.558321333333333
$ws.Range("N2").Value = 10.674964
$ws.Range("O2").Value = 0.3039644761000113
$ws.Range("P2").Value = 0.3039644761000113
$ws.Range("Q2").Value = 1.981884163562222
$ws.Range("R2").Value = 17.83695747206
$ws.Range("S2").Value = 0.104130008556077
$ws.Range("T2").Value = 0.104130008556077

$ws.Range("G3").Value = 0.5569716666666666
$ws.Range("H3").Value = 1.670915
$ws.Range("I3").Value = 0.3425729542218473
$ws.Range("J3").Value = 0.3425729542218473
$ws.Range("M3").Value = 5.383140666666667
$ws.Range("O3").Value = 0.4598470400038817
$ws.Range("P3").Value = 0.4598470400038817
$ws.Range("Q3").Value = 2.998256829014444
$ws.Range("R3").Value = 26.98431146113
$ws.Range("S3").Value = 0.1575311589843018
$ws.Range("T3").Value = 0.1575311589843017

$ws.Range("G4").Value = 0.5569716666666666
$ws.Range("H4").Value = 1.670915
$ws.Range("I4").Value = 0.3425729542218473
$ws.Range("J4").Value = 0.3425729542218473
$ws.Range("M4").Value = 2.764910333333333
$ws.Range("N4").Value = 8.294730999999999
$ws.Range("O4").Value = 0.2361884838961071
$ws.Range("P4").Value = 0.236188483896107
$ws.Range("Q4").Value = 1.539976716540555
$ws.Range("R4").Value = 13.859790448865
$ws.Range("S4").Value = 0.0809117866814686
$ws.Range("T4").Value = 0.08091178668146859

$ws.Range("G5").Value = 0.6936943333333333
$ws.Range("I5").Value = 0.4266660789393025
$ws.Range("J5").Value = 0.4266660789393025
$ws.Range("M5").Value = 3.558321333333333
$ws.Range("N5").Value = 10.674964
$ws.Range("O5").Value = 0.3039644761000113
$ws.Range("P5").Value = 0.3039644761000113
$ws.Range("Q5").Value = 2.468387345112444
$ws.Range("R5").Value = 22.215486106012
$ws.Range("S5").Value = 0.1296913311544312
$ws.Range("T5").Value = 0.1296913311544311

$ws.Range("G6").Value = 0.6936943333333333
$ws.Range("I6").Value = 0.4266660789393025
$ws.Range("J6").Value = 0.4266660789393025
$ws.Range("M6").Value = 5.383140666666667
$ws.Range("O6").Value = 0.4598470400038817
$ws.Range("P6").Value = 0.4598470400038817
$ws.Range("Q6").Value = 3.734254176002889
$ws.Range("R6").Value = 33.608287584026
$ws.Range("S6").Value = 0.1962011334703008
$ws.Range("T6").Value = 0.1962011334703007

$ws.Range("G7").Value = 0.6936943333333333
$ws.Range("I7").Value = 0.4266660789393025
$ws.Range("J7").Value = 0.4266660789393025
$ws.Range("M7").Value = 2.764910333333333
$ws.Range("N7").Value = 8.294730999999999
$ws.Range("O7").Value = 0.2361884838961071
$ws.Range("P7").Value = 0.236188483896107
$ws.Range("Q7").Value = 1.918002630408111
$ws.Range("R7").Value = 17.262023673673
$ws.Range("S7").Value = 0.1007736143145706
$ws.Range("T7").Value = 0.1007736143145706

$ws.Range("G8").Value = 0.3751823333333333
$ws.Range("H8").Value = 1.125547
$ws.Range("I8").Value = 0.2307609668388503
$ws.Range("J8").Value = 0.2307609668388503
$ws.Range("M8").Value = 3.558321333333333
$ws.Range("N8").Value = 10.674964
$ws.Range("O8").Value = 0.3039644761000113
$ws.Range("P8").Value = 0.3039644761000113
$ws.Range("Q8").Value = 1.335019300589778
$ws.Range("R8").Value = 12.015173705308
$ws.Range("S8").Value = 0.07014313638950322
$ws.Range("T8").Value = 0.0701431363895032

$ws.Range("G9").Value = 0.3751823333333333
$ws.Range("H9").Value = 1.125547
$ws.Range("I9").Value = 0.2307609668388503
$ws.Range("J9").Value = 0.2307609668388503
$ws.Range("M9").Value = 5.383140666666667
$ws.Range("O9").Value = 0.4598470400038817
$ws.Range("P9").Value = 0.4598470400038817
$ws.Range("Q9").Value = 2.019659275981556
$ws.Range("R9").Value = 18.176933483834
$ws.Range("S9").Value = 0.1061147475492792
$ws.Range("T9").Value = 0.1061147475492792

$ws.Range("G10").Value = 0.3751823333333333
$ws.Range("H10").Value = 1.125547
$ws.Range("I10").Value = 0.2307609668388503
$ws.Range("J10").Value = 0.2307609668388503
$ws.Range("M10").Value = 2.764910333333333
$ws.Range("N10").Value = 8.294730999999999
$ws.Range("O10").Value = 0.2361884838961071
$ws.Range("P10").Value = 0.236188483896107
$ws.Range("Q10").Value = 1.037345510317444
$ws.Range("R10").Value = 9.336109592857
$ws.Range("S10").Value = 0.05450308290006789
$ws.Range("T10").Value = 0.05450308290006788

